$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "og_group_ref"
$ws.Range("B17").Value = 144795
$ws.Range("D17").Value = $true

$ws.Range("D18").Select() | Out-Null
